$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values "following Dr Hou advice": ligand/receptor-expressing cell counts
# (columns E and K) change from 1 to 3 for all data rows, and the dependent
# expression/specificity statistics (G,H,I,J,M,N,Q,R,S,T) are recomputed
# accordingly for rows 2-5.

$data = @{
    2 = @{
        E = 3
        G = 2.943703
        H = 8.831109
        I = 0.1934541878053996
        J = 0.1934541878053996
        K = 3
        M = 2.761807333333334
        N = 8.285422000000001
        Q = 8.129940532555333
        R = 73.169464792998
        S = 0.1934541878053996
        T = 0.1934541878053996
    }
    3 = @{
        E = 3
        G = 2.508890333333333
        H = 7.526671
        I = 0.1648791816728176
        J = 0.1648791816728177
        K = 3
        M = 2.761807333333334
        N = 8.285422000000001
        Q = 6.929071721129112
        R = 62.361645490162
        S = 0.1648791816728176
        T = 0.1648791816728177
    }
    4 = @{
        E = 3
        G = 8.980481333333334
        H = 26.941444
        I = 0.5901790100569088
        J = 0.5901790100569088
        K = 3
        M = 2.761807333333334
        N = 8.285422000000001
        Q = 24.80235920326311
        R = 223.221232829368
        S = 0.5901790100569088
        T = 0.5901790100569088
    }
    5 = @{
        E = 3
        G = 0.7834633333333333
        H = 2.35039
        I = 0.05148762046487403
        J = 0.05148762046487404
        K = 3
        M = 2.761807333333334
        N = 8.285422000000001
        Q = 2.163774779397778
        R = 19.47397301458
        S = 0.05148762046487403
        T = 0.05148762046487404
    }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
